$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '36.700.37'
Set-TextValue 'E2' '  -0.97%  '
Set-TextValue 'D3' '2.059.90'
Set-TextValue 'E3' '  +0.57%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '243.85'
Set-TextValue 'E5' '  -0.93%  '
Set-TextValue 'E6' '  +1.26%  '
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'D8' '55.30'
Set-TextValue 'E8' '  -5.96%  '
Set-TextValue 'D9' '59.98'
Set-TextValue 'E9' '  +1.76%  '
Set-TextValue 'D10' '0.365'
Set-TextValue 'E10' '  -3.60%  '
Set-TextValue 'D11' '0.0752'
Set-TextValue 'E11' '  -2.79%  '
Set-TextValue 'E12' '  -3.10%  '
Set-TextValue 'D13' '0.935'
Set-TextValue 'E13' '  +4.79%  '
Set-TextValue 'D14' '14.79'
Set-TextValue 'E14' '  -3.94%  '
Set-TextValue 'D15' '2.358.17'
Set-TextValue 'E15' '  +0.74%  '
Set-TextValue 'D16' '5.46'
Set-TextValue 'E16' '  -4.65%  '
Set-TextValue 'D17' '2.068.37'
Set-TextValue 'E17' '  +1.87%  '
Set-TextValue 'D18' '36.612.88'
Set-TextValue 'E18' '  -1.15%  '
Set-TextValue 'E19' '  -6.69%  '
Set-TextValue 'D20' '72.17'
Set-TextValue 'E20' '  -2.24%  '
Set-TextValue 'D21' '0.0₃0864'
Set-TextValue 'E21' '  -2.45%  '
Set-TextValue 'D22' '238.43'
Set-TextValue 'E22' '  -0.82%  '
Set-TextValue 'D23' '5.27'
Set-TextValue 'E23' '  -2.51%  '
Set-TextValue 'E24' '  +0.04%  '
Set-TextValue 'E25' '  -2.90%  '
Set-TextValue 'D26' '2.14'
Set-TextValue 'E26' '  -0.39%  '
Set-TextValue 'E27' '  -3.48%  '
Set-TextValue 'D28' '165.00'
Set-TextValue 'E28' '  -2.11%  '
Set-TextValue 'D29' '20.19'
Set-TextValue 'E29' '  +0.92%  '
Set-TextValue 'E30' '  -1.64%  '
Set-TextValue 'E31' '  +6.88%  '
Set-TextValue 'D32' '5.11'
Set-TextValue 'E32' '  -7.82%  '
Set-TextValue 'D33' '4.52'
Set-TextValue 'E33' '  -4.79%  '
Set-TextValue 'D34' '0.0599'
Set-TextValue 'E34' '  -2.42%  '
Set-TextValue 'E35' '  -0.11%  '
Set-TextValue 'D36' '1.82'
Set-TextValue 'E36' '  -1.08%  '
Set-TextValue 'E37' '  -0.41%  '
Set-TextValue 'E38' '  -1.51%  '
Set-TextValue 'E39' '  -4.62%  '
Set-TextValue 'D40' '4.99'
Set-TextValue 'E40' '  -5.02%  '
Set-TextValue 'E41' '  -6.44%  '
Set-TextValue 'D42' '0.0216'
Set-TextValue 'E42' '  -3.10%  '
Set-TextValue 'E43' '  -3.27%  '
Set-TextValue 'D44' '94.73'
Set-TextValue 'E44' '  -2.87%  '
Set-TextValue 'B45' 'Maker'
Set-TextValue 'C45' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D45' '1.417.00'
Set-TextValue 'E45' '  +9.05%  '
Set-TextValue 'B46' 'Cronos'
Set-TextValue 'C46' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D46' '0.0910'
Set-TextValue 'E46' '  -6.29%  '
Set-TextValue 'D47' '7.68'
Set-TextValue 'E47' '  +13.80%  '
Set-TextValue 'D48' '16.14'
Set-TextValue 'E48' '  -5.17%  '
Set-TextValue 'E49' '  +1.59%  '
Set-TextValue 'E50' '  -4.32%  '
Set-TextValue 'D51' '2.248.17'
Set-TextValue 'E51' '  +0.98%  '
